$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cryptos list update (commit: Updated cryptos list on Fri Mar 10 20:42:41 UTC 2023 with GitHub Actions)

# Row 2
$ws.Range("D2").Value = "'19.972.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.75%  "

# Row 3
$ws.Range("D3").Value = "'1.418.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.81%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.82%  "

# Row 5
$ws.Range("D5").Value = "'1.000"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.68%  "

# Row 6
$ws.Range("D6").Value = "'276.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.09%  "

# Row 7
$ws.Range("D7").Value = "'0.3689"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.93%  "

# Row 8
$ws.Range("D8").Value = "'0.3105"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.76%  "

# Row 9
$ws.Range("D9").Value = "'39.83"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.30%  "

# Row 10
$ws.Range("D10").Value = "'1.042"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.73%  "

# Row 11
$ws.Range("D11").Value = "'0.06544"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.79%  "

# Row 12
$ws.Range("D12").Value = "'0.9998"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.01%  "

# Row 13
$ws.Range("D13").Value = "'5.507"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.56%  "

# Row 14
$ws.Range("D14").Value = "'17.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.50%  "

# Row 15
$ws.Range("D15").Value = "'6.221"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.41%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'1.422.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.35%  "

# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.00001023"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.94%  "

# Row 18
$ws.Range("D18").Value = "'0.05692"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -12.21%  "

# Row 19
$ws.Range("D19").Value = "'0.9997"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.73%  "

# Row 20
$ws.Range("D20").Value = "'71.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -10.66%  "

# Row 21
$ws.Range("D21").Value = "'5.623"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.98%  "

# Row 22
$ws.Range("D22").Value = "'14.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "

# Row 23
$ws.Range("D23").Value = "'10.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.20%  "

# Row 24
$ws.Range("D24").Value = "'2.239"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.85%  "

# Row 25
$ws.Range("D25").Value = "'20.012.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.57%  "

# Row 26
$ws.Range("D26").Value = "'2.287"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.49%  "

# Row 27
$ws.Range("D27").Value = "'133.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.23%  "

# Row 28
$ws.Range("D28").Value = "'17.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.44%  "

# Row 29
$ws.Range("D29").Value = "'1.581.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.59%  "

# Row 30
$ws.Range("D30").Value = "'110.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.51%  "

# Row 31
$ws.Range("D31").Value = "'3.897"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -18.38%  "

# Row 32
$ws.Range("D32").Value = "'5.246"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.04%  "

# Row 33
$ws.Range("D33").Value = "'0.8213"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.87%  "

# Row 34
$ws.Range("D34").Value = "'0.07767"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.18%  "

# Row 35
$ws.Range("D35").Value = "'1.485"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.08%  "

# Row 36
$ws.Range("D36").Value = "'8.201"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.04%  "

# Row 37
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.05870"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.49%  "

# Row 38
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'4.903"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.87%  "

# Row 40
$ws.Range("D40").Value = "'0.02057"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.50%  "

# Row 41
$ws.Range("D41").Value = "'10.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.48%  "

# Row 42
$ws.Range("D42").Value = "'0.1885"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.97%  "

# Row 43
$ws.Range("D43").Value = "'1.103"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.62%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'12.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.24%  "

# Row 45
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.5327"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.40%  "

# Row 46
$ws.Range("D46").Value = "'3.540"
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'116.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.30%  "

# Row 48
$ws.Range("D48").Value = "'0.5212"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.17%  "

# Row 49
$ws.Range("D49").Value = "'1.773"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.00%  "

# Row 50
$ws.Range("D50").Value = "'1.035"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.32%  "

# Row 51
$ws.Range("D51").Value = "'0.9997"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.71%  "
